# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The account-statement database (columns E "Periodo Mora" and F on rows
# 16-23) is re-sorted: the oldest period (2012) that used to sit on top
# (row 16) moves to the bottom (row 23), and the newest period (2107)
# that used to be last moves to the top - i.e. the eight existing rows
# are reversed in place, carrying their "fecha/valor" (F column) along
# with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2107", "2106", "2105", "2104", "2103", "2102", "2101", "2012")
$fvals   = @(29260, 35112, 35112, 35112, 35112, 35112, 35112, 35112)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $fvals[$i]
}

# The wider "2107"/"2106"/... column and the new date/amount values shift
# the best-fit column widths very slightly; re-apply the closest widths.
$ws.Columns.Item(2).ColumnWidth  = 17.6666666666667
$ws.Columns.Item(3).ColumnWidth  = 15.8333333333333
$ws.Columns.Item(5).ColumnWidth  = 12.6666666666667
$ws.Columns.Item(6).ColumnWidth  = 9.33333333333333
$ws.Columns.Item(7).ColumnWidth  = 13.5
$ws.Columns.Item(8).ColumnWidth  = 18.5
$ws.Columns.Item(9).ColumnWidth  = 17.3333333333333
$ws.Columns.Item(10).ColumnWidth = 14.1666666666667

$wb.Save()
